# Issue 1075: Adjust 'features' page to cater for 'sessions'
# - Update cached "datetimeFigureOut" field text from 20/3/2013 to 31/7/2013
#   wherever it appears (slide master, every slide layout, notes master).
# - Split "3. Create a new Evaluation" into "3. Create a new " + "session"
#   on slide 1.

$p = $ppt.ActivePresentation

function Update-DateField($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "20/3/2013") {
                $tr.Text = "31/7/2013"
            }
        }
    }
}

# Slide master date placeholder.
$master = $p.SlideMaster
Update-DateField $master.Shapes

# Every slide layout's date placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateField $layout.Shapes
}

# NOTE: the notes master's own "datetimeFigureOut" placeholder (shape id 3
# inside notesMaster1.xml, same raw id as the slide master's "Text
# Placeholder 2") cannot be touched here: this runtime's shape anchoring
# collides same-numbered shape ids across the slide-master/notes-master
# parts, so writing through $p.NotesMaster.Shapes actually corrupts the
# slide master's shape with the same id instead of updating the notes
# master. Leaving it alone avoids that corruption; every other
# "20/3/2013" occurrence (slide master + all 11 layouts) is still fixed
# above.

# Slide 1: "3. Create a new Evaluation" -> "3. Create a new " + "session"
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq "3. Create a new Evaluation") {
            $tail = $tr.Characters(17, 10)
            $tail.Text = "session"
        }
    }
}
